{"js": "// Replace 25 division-problem answers in the table, per the diff.\n// Strategy: locate every target Range FIRST (search is scoped against the\n// original, pre-edit text), then apply all insertText replacements using\n// those captured Range objects. This avoids any cross-talk between pairs\n// where a pair's new value happens to equal another pair's old value\n// (e.g. '56\u00f73=18, 2' -> '62\u00f73=20, 2' while a different cell already\n// holds '62\u00f73=20, 2' and must become '66\u00f78=8, 2').\nconst replacements = [\n  { oldText: \"95\u00f79=10, 5\", newText: \"66\u00f75=13, 1\" },\n  { oldText: \"56\u00f73=18, 2\", newText: \"62\u00f73=20, 2\" },\n  { oldText: \"43\u00f79=4, 7\", newText: \"57\u00f73=19, 0\" },\n  { oldText: \"74\u00f78=9, 2\", newText: \"87\u00f73=29, 0\" },\n  { oldText: \"41\u00f78=5, 1\", newText: \"23\u00f76=3, 5\" },\n  { oldText: \"50\u00f75=10, 0\", newText: \"66\u00f73=22, 0\" },\n  { oldText: \"84\u00f73=28, 0\", newText: \"50\u00f74=12, 2\" },\n  { oldText: \"25\u00f77=3, 4\", newText: \"43\u00f73=14, 1\" },\n  { oldText: \"87\u00f74=21, 3\", newText: \"43\u00f79=4, 7\" },\n  { oldText: \"36\u00f75=7, 1\", newText: \"33\u00f73=11, 0\" },\n  { oldText: \"28\u00f74=7, 0\", newText: \"15\u00f75=3, 0\" },\n  { oldText: \"82\u00f76=13, 4\", newText: \"99\u00f72=49, 1\" },\n  { oldText: \"70\u00f74=17, 2\", newText: \"78\u00f72=39, 0\" },\n  { oldText: \"62\u00f73=20, 2\", newText: \"66\u00f78=8, 2\" },\n  { oldText: \"36\u00f77=5, 1\", newText: \"88\u00f74=22, 0\" },\n  { oldText: \"57\u00f75=11, 2\", newText: \"96\u00f72=48, 0\" },\n  { oldText: \"88\u00f74=22, 0\", newText: \"28\u00f78=3, 4\" },\n  { oldText: \"74\u00f72=37, 0\", newText: \"73\u00f77=10, 3\" },\n  { oldText: \"68\u00f77=9, 5\", newText: \"81\u00f73=27, 0\" },\n  { oldText: \"90\u00f72=45, 0\", newText: \"75\u00f73=25, 0\" },\n  { oldText: \"26\u00f76=4, 2\", newText: \"33\u00f77=4, 5\" },\n  { oldText: \"28\u00f76=4, 4\", newText: \"54\u00f79=6, 0\" },\n  { oldText: \"96\u00f73=32, 0\", newText: \"69\u00f76=11, 3\" },\n  { oldText: \"44\u00f79=4, 8\", newText: \"47\u00f76=7, 5\" },\n  { oldText: \"41\u00f74=10, 1\", newText: \"52\u00f76=8, 4\" },\n];\n\n// Step 1: search for each old text (unique, single match each) and keep\n// the resulting Range objects around.\nconst searchResults = replacements.map(r =>\n  context.document.body.search(r.oldText, { matchCase: true, ignorePunct: false })\n);\nsearchResults.forEach(res => res.load('items'));\nawait context.sync();\n\n// Step 2: verify each search found exactly one match, then replace its text.\nsearchResults.forEach((res, i) => {\n  if (res.items.length !== 1) {\n    throw new Error(\n      `Expected exactly 1 match for \"${replacements[i].oldText}\", found ${res.items.length}`\n    );\n  }\n  res.items[0].insertText(replacements[i].newText, Word.InsertLocation.replace);\n});\nawait context.sync();\n", "ps1": "# Update the 25 division-problem answers in the table (5 answer-rows x 5 columns),\n# per the diff. Each entry is addressed by its fixed (row, column) position in the\n# table rather than by searching for the old text, because several of the new\n# values collide with other cells' old values (e.g. one cell goes\n# 56\u00f73=18, 2 -> 62\u00f73=20, 2 while another cell already holds 62\u00f73=20, 2 and\n# must become 66\u00f78=8, 2). Position-based addressing sidesteps that ambiguity\n# entirely.\n\n$d = $word.ActiveDocument\n$tbl = $d.Tables.Item(1)\n\n# row, column, expected old text, new text\n$updates = @(\n    @(1, 1, \"95\u00f79=10, 5\", \"66\u00f75=13, 1\"),\n    @(1, 2, \"56\u00f73=18, 2\", \"62\u00f73=20, 2\"),\n    @(1, 3, \"43\u00f79=4, 7\", \"57\u00f73=19, 0\"),\n    @(1, 4, \"74\u00f78=9, 2\", \"87\u00f73=29, 0\"),\n    @(1, 5, \"41\u00f78=5, 1\", \"23\u00f76=3, 5\"),\n    @(5, 1, \"50\u00f75=10, 0\", \"66\u00f73=22, 0\"),\n    @(5, 2, \"84\u00f73=28, 0\", \"50\u00f74=12, 2\"),\n    @(5, 3, \"25\u00f77=3, 4\", \"43\u00f73=14, 1\"),\n    @(5, 4, \"87\u00f74=21, 3\", \"43\u00f79=4, 7\"),\n    @(5, 5, \"36\u00f75=7, 1\", \"33\u00f73=11, 0\"),\n    @(9, 1, \"28\u00f74=7, 0\", \"15\u00f75=3, 0\"),\n    @(9, 2, \"82\u00f76=13, 4\", \"99\u00f72=49, 1\"),\n    @(9, 3, \"70\u00f74=17, 2\", \"78\u00f72=39, 0\"),\n    @(9, 4, \"62\u00f73=20, 2\", \"66\u00f78=8, 2\"),\n    @(9, 5, \"36\u00f77=5, 1\", \"88\u00f74=22, 0\"),\n    @(13, 1, \"57\u00f75=11, 2\", \"96\u00f72=48, 0\"),\n    @(13, 2, \"88\u00f74=22, 0\", \"28\u00f78=3, 4\"),\n    @(13, 3, \"74\u00f72=37, 0\", \"73\u00f77=10, 3\"),\n    @(13, 4, \"68\u00f77=9, 5\", \"81\u00f73=27, 0\"),\n    @(13, 5, \"90\u00f72=45, 0\", \"75\u00f73=25, 0\"),\n    @(17, 1, \"26\u00f76=4, 2\", \"33\u00f77=4, 5\"),\n    @(17, 2, \"28\u00f76=4, 4\", \"54\u00f79=6, 0\"),\n    @(17, 3, \"96\u00f73=32, 0\", \"69\u00f76=11, 3\"),\n    @(17, 4, \"44\u00f79=4, 8\", \"47\u00f76=7, 5\"),\n    @(17, 5, \"41\u00f74=10, 1\", \"52\u00f76=8, 4\"),\n)\n\nforeach ($u in $updates) {\n    $row = $u[0]\n    $col = $u[1]\n    $oldText = $u[2]\n    $newText = $u[3]\n\n    $cell = $tbl.Cell($row, $col)\n    $current = $cell.Range.Text.TrimEnd([char]13, [char]7)\n    if ($current -ne $oldText) {\n        throw \"Cell ($row,$col): expected '$oldText' but found '$current'\"\n    }\n    $cell.Range.Text = $newText\n}\n\nWrite-Output \"Replaced $($updates.Count) cells.\"\n"}
